$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Summary" to "Sheet1"
$ws.Name = "Sheet1"

# Update Universe_Returns column (C) for rows 2-6 with higher-precision value
$ws.Range("C2:C6").Value = 6.039999961853027

# Remove rows 7-11 (the extra strategy rows no longer part of the summary)
$ws.Range("A7:K11").EntireRow.Delete()
